# "Added Games to Presentation"
#
# The original deck has a single Title-Slide-layout slide (ctrTitle +
# subTitle placeholders, both empty). The edit turns it into five
# "Title Only" slides, each with just a title placeholder naming a game:
#   1. Elden Ring
#   2. Dark Souls 3
#   3. Red Dead Redemption 2
#   4. RimWorld
#   5. Minecraft

$p = $ppt.ActivePresentation

# --- Slide 1: swap the Title-Slide layout content for a Title-Only one ---
# Changing .Layout / .CustomLayout in place on the existing slide does not
# rebuild its placeholders, so instead we insert a fresh "Title Only" slide
# (ppLayoutTitleOnly = 11) - which is born with just the single title
# placeholder we want - move it to the front, and drop the original
# ctrTitle/subTitle slide.
$freshSlide = $p.Slides.Add(2, 11)
$freshSlide.MoveTo(1)
$p.Slides.Item(2).Delete()

$slide1 = $p.Slides.Item(1)
$slide1.Shapes.Item(1).TextFrame.TextRange.Text = "Elden Ring"

# --- Slides 2-5: duplicate slide 1 for each remaining game ---
# Each Duplicate() call inserts its copy right after slide 1, so duplicating
# slide 1 repeatedly (instead of duplicating the newest slide) stacks the
# newest duplicate closest to slide 1; queue the titles last-to-first so the
# final left-to-right order reads correctly.
$titlesNewestFirst = @("Minecraft", "RimWorld", "Red Dead Redemption 2", "Dark Souls 3")

foreach ($title in $titlesNewestFirst) {
    $dup = $slide1.Duplicate()
    $dup.Item(1).Shapes.Item(1).TextFrame.TextRange.Text = $title
}
